$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Jumper row (row 12): quantity 1 -> 3, designator J1 -> J1,J2,J3
$ws.Range("B12").Value = 3
$ws.Range("D12").Value = "J1,J2,J3"

# 2) Remove the R4,R5 22-ohm resistor row (row 20) entirely -- no longer
#    required now that track width / Cu weight (2oz/ft^2) handles the
#    current budget without series limiting resistors.
$ws.Rows("20:20").Delete()

# 3) Renumber the "No." column (A) for all rows that shifted up, so the
#    sequential numbering stays contiguous (1..32) after the deletion.
#    Row 33 is a trailing blank-string row (not a real BOM line) and must
#    stay untouched.
for ($r = 20; $r -le 32; $r++) {
    $ws.Cells.Item($r, 1).Value = ($r - 1).ToString()
}
